$d = $word.ActiveDocument

# Step 1: Replace the text of the first run ("EM KHÔNG BIẾ") with the full
# new sentence.
$r1 = $d.Content
$found1 = $r1.Find.Execute("EM KHÔNG BIẾ", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "EM THƯA THẦY MẤY BÀI NÀY EM CHƯA BIẾT VẼ Ạ (CÓ GÌ EM LÊN LỚP THẦY CHỈ BAAOR EM THÊM Ạ)", `
    2)

# Step 2: Remove the trailing run that used to read "T VẼ Ạ". It sits right
# after the bookmark, at the very end of the paragraph/document, so build an
# explicit range from the end of the just-replaced text through the end of
# the document content and delete it outright (this removes the whole run,
# leaving the bookmark untouched).
if ($found1) {
    $tail = $d.Range($r1.End, $d.Content.End)
    if ($tail.Start -lt $tail.End) {
        $tail.Delete()
    }
}
